$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.47
$ws.Range("B3").Value = 7.502000000000001
$ws.Range("B4").Value = 5.3245
$ws.Range("B5").Value = 7.18
